# Insert a new weekly record at row 42 ("Fruta / hortaliza, semanal").
# This shifts the existing rows 42..136 down to 43..137 and the new
# row 42 receives a fresh set of values while keeping every other
# (constant) column identical to the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 42..136 down to 43..137, leaving row 42 blank.
$ws.Rows.Item(42).Insert()

# Fill the new row 42 with the new weekly record.
$ws.Cells.Item(42, 1).Value  = 8
$ws.Cells.Item(42, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(42, 3).Value  = "Coquimbo"
$ws.Cells.Item(42, 4).Value  = 45082
$ws.Cells.Item(42, 5).Value  = 4
$ws.Cells.Item(42, 6).Value  = 100114007
$ws.Cells.Item(42, 7).Value  = "Jengibre"
$ws.Cells.Item(42, 8).Value  = "Sin especificar"
$ws.Cells.Item(42, 9).Value  = "Primera"
$ws.Cells.Item(42, 10).Value = 440
$ws.Cells.Item(42, 11).Value = 17000
$ws.Cells.Item(42, 12).Value = 18000
$ws.Cells.Item(42, 13).Value = 17500
$ws.Cells.Item(42, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(42, 15).Value = "Perú"
$ws.Cells.Item(42, 16).Value = 1346
$ws.Cells.Item(42, 17).Value = 13
$ws.Cells.Item(42, 18).Value = "Hortaliza"

# Match the date-number format used by the rest of the "Fecha" column.
$ws.Cells.Item(42, 4).NumberFormat = $ws.Cells.Item(43, 4).NumberFormat
